# hw12-ChristianHansen.xlsx — "Add files via upload" re-save.
#
# Net effect versus the previous version (per the OOXML diff):
#   - B4 (Runtime for "Localizing parameters in insertNode") changes from
#     17.209734000000001 to 28.084282999999999.
#   - C4 ("<- holy" annotation / callout next to that data point) is removed,
#     which also drops "<- holy" from the shared-strings table and shrinks
#     the sheet's used range from A1:C8 down to A1:B8.
#   - The sheet's active-cell selection ends up on C4 (the now-empty cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Correct the runtime figure for the "Localizing parameters in insertNode" row.
$ws.Range("B4").Value = 28.084282999999999

# Drop the "<- holy" callout entirely — the cell goes back to empty, which
# also removes the now-unused shared string and shrinks the used range.
$ws.Range("C4").ClearContents()

# Leave the selection sitting on the now-empty C4, matching the saved file.
$ws.Range("C4").Select()
